$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently follows the title.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. Insert a new bold "Play African Quest for Free - A Detailed Review" paragraph
#    right before the final (DALL-E prompt) paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play African Quest for Free ' + [char]0x2013 + ' A Detailed Review</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$insertPoint.InsertXML($titleXml) | Out-Null

# InsertXML left a trailing empty placeholder paragraph between the new title
# paragraph and the original last paragraph - remove it.
$placeholder = $d.Paragraphs($count + 1)
$placeholder.Range.Delete()

# 3. Replace the DALL-E image-prompt text with the new meta-description copy,
#    keeping the run's existing (italic) formatting.
$finalCount = $d.Paragraphs.Count
$dallePara = $d.Paragraphs($finalCount)
$pr = $dallePara.Range
$target = $d.Range($pr.Start, $pr.End)
$target.Text = "Discover African Quest Slot Machine's pros and cons, volatility, RTP, special features, and target audience. Play for free at recommended casinos."
